$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 10, shifting existing rows 10..84 down to 11..85
$ws.Rows("10:10").Insert()

# Copy the "constant" columns from the (now shifted) row below so the new
# row matches the rest of the dataset, then set the row-specific values.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Terminal La Palmera de La Serena"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44473
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112001
$ws.Range("G10").Value = "Berenjena"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 700
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7500
$ws.Range("N10").Value = "$/caja 60 unidades"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 125
$ws.Range("Q10").Value = 60
$ws.Range("R10").Value = "Hortaliza"
